$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "AVG NUMBER OF WORDS PER SENTENCE" column (column H).
# Deleting the entire column shifts I:M left to H:L, matching the diff.
$ws.Range("H1").EntireColumn.Delete()
